$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns remain text, so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "51.584.84"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.791.85"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "351.17"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "108.69"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  +5.48%  "
$ws.Range("D10").Value = "39.57"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "19.98"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "3.226.58"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "2.794.22"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "0.932"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "51.563.96"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("D20").Value = "3.13"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "13.39"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").Value = "70.51"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "267.01"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "2.74"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "25.88"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "10.30"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "37.30"
$ws.Range("E30").Value = "  +8.50%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("D33").Value = "51.97"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "5.68"
$ws.Range("E34").Value = "  +8.59%  "
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "18.59"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").Value = "1.97"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "2.49"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").Value = "120.12"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "2.19"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "21.90"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "2.133.28"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").Value = "3.36"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  +5.67%  "
$ws.Range("D49").Value = "0.224"
$ws.Range("E49").Value = "  +17.60%  "
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "0.901"
$ws.Range("E50").Value = "  -5.17%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "1.35"
$ws.Range("E51").Value = "  +8.93%  "
